$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.283.27'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '3.413.84'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '254.89'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '671.71'
$ws.Range('E6').Value = '  -1.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.52'
$ws.Range('E7').Value = '  +5.41%  '
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.07'
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').Value = '3.411.84'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.90'
$ws.Range('E12').Value = '  +10.47%  '
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').Value = '98.077.32'
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('E15').Value = '  -1.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000260'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '4.049.87'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('E18').Value = '  +3.54%  '
$ws.Range('D19').Value = '3.410.36'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.48'
$ws.Range('E20').Value = '  +6.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.542'
$ws.Range('E21').Value = '  -4.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.45'
$ws.Range('E22').Value = '  +4.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '514.05'
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.44'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000203'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  +5.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '97.94'
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.53'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').Value = '3.588.89'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.35'
$ws.Range('E30').Value = '  +7.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.95'
$ws.Range('E31').Value = '  +12.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.145'
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  -2.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.571'
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '29.30'
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.02'
$ws.Range('E38').Value = '  +1.61%  '
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '531.13'
$ws.Range('E40').Value = '  +1.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.155'
$ws.Range('E41').Value = '  +1.76%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.869'
$ws.Range('E43').Value = '  +0.71%  '
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.77'
$ws.Range('E45').Value = '  +2.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0429'
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.69'
$ws.Range('E47').Value = '  -2.87%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.68'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.70'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '56.11'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.23'
$ws.Range('E51').Value = '  +6.26%  '
